# Generate Report for Handoff
#
# The localization run regenerated the handoff package: the e2e test
# document got a new GUID-based file name, and new handoff timestamps /
# xliff file names were produced. Propagate the new values into the
# status report workbook (Overview + per-locale sheets), including the
# hyperlinks that point at the e2e markdown file.

$wb = $excel.ActiveWorkbook

$oldGuid = "1b4b7c88-8ffc-4aa3-a3ab-2ebe1ff954cb"
$newGuid = "bcbcd2ce-0db8-4db9-ae74-df2c4b3632dd"

$oldMd   = "$oldGuid.md"
$newMd   = "$newGuid.md"

$oldHash = "4ca0f1ee8dc08a8bc851482dd8d18d01ecb9820a"
$newHash = "20c32e1f0550ab60d279708d1fdc422defca2bb6"

# The external hyperlink target (points at the historical git blob) is not
# part of this change - it stays addressed at the old file name.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/974d3b448c736b5329d1a16ceb0f03e6b16b741d/e2e/$oldMd"

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# ---- Overview sheet ----------------------------------------------------
# A2: File Name
$overview.Range("A2").Value = $newMd

# B2: Path And Name (value + hyperlink display text)
$overview.Range("B2").Value = "e2e\$newMd"
$overview.Range("B2").Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("B2"), $linkAddress, "", "", "e2e\$newMd")

# G2: Latest HO Xliff Generate Date
$overview.Range("G2").Value = "2016-08-26 19:00:23"

# ---- zh-cn sheet ---------------------------------------------------------
# A2: Source File Name (value + hyperlink display text)
$zhcn.Range("A2").Value = $newMd
$zhcn.Range("A2").Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $linkAddress, "", "", $newMd)

# G2: Latest Handoff File
$zhcn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"

# H2: Latest Handoff Datetime
$zhcn.Range("H2").Value = "2016-08-26 19:00:00"

# ---- de-de sheet ---------------------------------------------------------
# A2: Source File Name (value + hyperlink display text)
$dede.Range("A2").Value = $newMd
$dede.Range("A2").Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $linkAddress, "", "", $newMd)

# G2: Latest Handback File
$dede.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"

# H2: Latest HO Xliff Generate Date (shared with Overview!G2 text)
$dede.Range("H2").Value = "2016-08-26 19:00:23"
